$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("B15").Value = "1-2 Weeks"
$ws.Range("C15").Value = "tiwariravikant04@gmail.com"
$ws.Range("D15").Value = "Ravikant Tiwari"
$ws.Range("E15").Value = "'8744883594"
$ws.Range("F15").Value = "'+40"
$ws.Range("G15").Value = "11/17/2024, 1:23:16 PM"

# Row 16
$ws.Range("A16").Value = "Life"
$ws.Range("B16").Value = "1-2 Weeks"
$ws.Range("C16").Value = "tiwariravikant04@gmail.com"
$ws.Range("D16").Value = "Ravikant Tiwari"
$ws.Range("E16").Value = "'8744883594"
$ws.Range("F16").Value = "'+40"
$ws.Range("G16").Value = "11/17/2024, 1:40:24 PM"

# Row 17
$ws.Range("B17").Value = "1week"
$ws.Range("C17").Value = "tiwariravikant04@gmail.com"
$ws.Range("D17").Value = "Ravikant Tiwari"
$ws.Range("E17").Value = "'8744883594"
$ws.Range("F17").Value = "'+40"

# Row 18
$ws.Range("A18").Value = "Travel"
$ws.Range("B18").Value = "1-2 Weeks"
$ws.Range("C18").Value = "ravikanttiwari488@gmail.com"
$ws.Range("D18").Value = "Ravikant Tiwari"
$ws.Range("E18").Value = "'8744883594"
$ws.Range("F18").Value = "'+40"
$ws.Range("G18").Value = "11/17/2024, 2:32:34 PM"
